$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 144, pushing existing rows 144-181 down to 145-182.
$ws.Rows.Item(144).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A144").Value = 4
$ws.Range("B144").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C144").Value = "Los Lagos"
$ws.Range("D144").Value = 45244
$ws.Range("E144").Value = 10
$ws.Range("F144").Value = 100112022
$ws.Range("G144").Value = "Arveja Verde"
$ws.Range("H144").Value = "Sin especificar"
$ws.Range("I144").Value = "Primera"
$ws.Range("J144").Value = 90
$ws.Range("K144").Value = 29000
$ws.Range("L144").Value = 29000
$ws.Range("M144").Value = 29000
$ws.Range("N144").Value = "$/saco 25 kilos"
$ws.Range("O144").Value = "Región del Maule"
$ws.Range("P144").Value = 1160
$ws.Range("Q144").Value = 25
$ws.Range("R144").Value = "Hortaliza"
